$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultiLevel")

$ws.Range("B5").Value = '<jt:if test="true">$[COUNTA(B3||$Z$1)]'
$ws.Range("F5").Value = '$[SUM(C3)/SUM(E3||1)]</jt:if>'
$ws.Range("I1").Value = '$[COUNTA(''Formula Test''!$E$3)]'
$ws.Range("I2").Value = '$[COUNTA(''Formula Test''!$K$3)]'
